$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.460.68'
$ws.Range('E2').Value = '  -2.71%  '
$ws.Range('D3').Value = '3.689.30'
$ws.Range('E3').Value = '  -3.23%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '687.88'
$ws.Range('E5').Value = '  -1.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.74'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('D7').Value = '3.690.23'
$ws.Range('E7').Value = '  -3.16%  '
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  -5.73%  '
$ws.Range('E10').Value = '  -8.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.35'
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.437'
$ws.Range('E12').Value = '  -9.88%  '
$ws.Range('E13').Value = '  -6.17%  '
$ws.Range('D14').Value = '4.311.24'
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.93'
$ws.Range('E15').Value = '  -8.68%  '
$ws.Range('D16').Value = '3.687.51'
$ws.Range('E16').Value = '  -3.39%  '
$ws.Range('D17').Value = '69.509.76'
$ws.Range('E17').Value = '  -2.63%  '
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.00'
$ws.Range('E19').Value = '  -8.90%  '
$ws.Range('E20').Value = '  -10.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '475.20'
$ws.Range('E21').Value = '  -7.71%  '
$ws.Range('E22').Value = '  -5.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.655'
$ws.Range('E23').Value = '  -8.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.84'
$ws.Range('E24').Value = '  -5.08%  '
$ws.Range('D25').Value = '3.834.90'
$ws.Range('E25').Value = '  -3.14%  '
$ws.Range('E26').Value = '  -9.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.19'
$ws.Range('E28').Value = '  -12.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.23'
$ws.Range('E29').Value = '  -10.92%  '
$ws.Range('E30').Value = '  -11.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.70'
$ws.Range('E31').Value = '  -10.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.75'
$ws.Range('E32').Value = '  -8.02%  '
$ws.Range('E33').Value = '  -8.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.78'
$ws.Range('E35').Value = '  -8.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.164'
$ws.Range('E36').Value = '  -5.02%  '
$ws.Range('D37').Value = '3.655.53'
$ws.Range('E37').Value = '  -3.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.30'
$ws.Range('E38').Value = '  -10.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.22'
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.30'
$ws.Range('E40').Value = '  -6.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0915'
$ws.Range('E42').Value = '  -9.59%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  -6.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '163.96'
$ws.Range('E45').Value = '  -5.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '48.30'
$ws.Range('E46').Value = '  -3.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '30.18'
$ws.Range('E47').Value = '  +2.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.74'
$ws.Range('E48').Value = '  -15.98%  '
$ws.Range('E49').Value = '  -4.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000279'
$ws.Range('E50').Value = '  -9.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.11'
$ws.Range('E51').Value = '  -4.55%  '

$wb.Save()
